$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.882.51'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.57%  '

# Row 3
$ws.Range("E3").Value = '  -2.29%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '599.37'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.41%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.63'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.35%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.860.18'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.14%  '

# Row 8
$ws.Range("E8").Value = '  +0.12%  '

# Row 9
$ws.Range("E9").Value = '  -1.43%  '

# Row 10
$ws.Range("E10").Value = '  -5.07%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.43'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.62%  '

# Row 12
$ws.Range("E12").Value = '  -2.63%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000260'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.88%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.91'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.66%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.504.48'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.11%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.851.93'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.40%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.060.57'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.16%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.97'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.10%  '

# Row 19
$ws.Range("E19").Value = '  -2.86%  '

# Row 20
$ws.Range("E20").Value = '  -0.62%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.76'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.84%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '465.73'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.89%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.732'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.52%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000158'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.16%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.87'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.97%  '

# Row 26
$ws.Range("E26").Value = '  -2.92%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.02'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.92%  '

# Row 28
$ws.Range("E28").Value = '  -0.20%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.92'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.47%  '

# Row 30
$ws.Range("E30").Value = '  -1.04%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.006.70'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.19%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.65'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.37%  '

# Row 33
$ws.Range("E33").Value = '  -4.74%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '31.09'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.45%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.46'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.81%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.821.02'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.55%  '

# Row 37
$ws.Range("E37").Value = '  -3.40%  '

# Row 38
$ws.Range("E38").Value = '  +11.04%  '

# Row 39
$ws.Range("E39").Value = '  -2.11%  '

# Row 40
$ws.Range("E40").Value = '  -0.16%  '

# Row 41
$ws.Range("E41").Value = '  -3.32%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.24%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.312'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.71%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '422.22'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.82%  '

# Row 45
$ws.Range("B45").Value = 'FLOKI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.000297'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.20%  '

# Row 46
$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.96'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.94%  '

# Row 47
$ws.Range("E47").Value = '  -0.02%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.59'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.27%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '46.98'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.64%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '26.59'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.48%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '142.05'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.66%  '
